$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.675.64"
$ws.Range("E2").Value = "  -1.80%  "
$ws.Range("D3").Value = "3.512.23"
$ws.Range("E3").Value = "  -2.56%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'586.59"
$ws.Range("E5").Value = "  -3.11%  "
$ws.Range("D6").Value = "'132.50"
$ws.Range("E6").Value = "  -3.05%  "
$ws.Range("D7").Value = "3.512.67"
$ws.Range("E7").Value = "  -2.46%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("E10").Value = "  -0.91%  "
$ws.Range("D11").Value = "'7.30"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("D13").Value = "4.112.41"
$ws.Range("E13").Value = "  -2.49%  "
$ws.Range("D14").Value = "'27.80"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("E15").Value = "  -3.76%  "
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "3.512.08"
$ws.Range("E17").Value = "  -2.43%  "
$ws.Range("D18").Value = "64.671.21"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").Value = "'9.96"
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("D20").Value = "'14.32"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("D21").Value = "'5.68"
$ws.Range("E21").Value = "  -4.54%  "
$ws.Range("D22").Value = "'392.71"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("E23").Value = "  -1.82%  "
$ws.Range("D24").Value = "3.654.20"
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("D25").Value = "'73.69"
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "'0.0000110"
$ws.Range("E27").Value = "  -6.51%  "
$ws.Range("E28").Value = "  -5.83%  "
$ws.Range("D29").Value = "'7.47"
$ws.Range("E29").Value = "  -7.94%  "
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("E31").Value = "  -5.69%  "
$ws.Range("D32").Value = "'8.21"
$ws.Range("E32").Value = "  -5.88%  "
$ws.Range("D33").Value = "3.513.21"
$ws.Range("E33").Value = "  -2.46%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "'24.07"
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("D37").Value = "'5.28"
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("D39").Value = "'171.02"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").Value = "'7.00"
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("E41").Value = "  -3.23%  "
$ws.Range("D42").Value = "'26.61"
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("E43").Value = "  -3.60%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "'42.14"
$ws.Range("E45").Value = "  -2.94%  "
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("E47").Value = "  -3.16%  "
$ws.Range("E48").Value = "  -2.98%  "
$ws.Range("D49").Value = "2.447.47"
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("E50").Value = "  -2.73%  "
$ws.Range("D51").Value = "'0.902"
$ws.Range("E51").Value = "  +2.00%  "
